$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted above the current
# first row of data for this product/market (row 367), pushing the
# existing rows (367-389) down by one, so the table grows from
# A1:R389 to A1:R390.
$ws.Rows.Item(367).Insert()

# Populate the newly inserted row 367 with the new observation.
$ws.Cells.Item(367, 1).Value = 9
$ws.Cells.Item(367, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(367, 3).Value = "Metropolitana"
$ws.Cells.Item(367, 4).Value = 44714
$ws.Cells.Item(367, 5).Value = 13
$ws.Cells.Item(367, 6).Value = 100112032
$ws.Cells.Item(367, 7).Value = "Zapallo italiano"
$ws.Cells.Item(367, 8).Value = "Sin especificar"
$ws.Cells.Item(367, 9).Value = "Primera"
$ws.Cells.Item(367, 10).Value = 160
$ws.Cells.Item(367, 11).Value = 13000
$ws.Cells.Item(367, 12).Value = 14000
$ws.Cells.Item(367, 13).Value = 13500
$ws.Cells.Item(367, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(367, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(367, 16).Value = 225
$ws.Cells.Item(367, 17).Value = 60
$ws.Cells.Item(367, 18).Value = "Hortaliza"
